$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows before the current last row (16) ---
# This pushes the existing row 16 ("Ajout d'un systeme de log") down to row 18,
# and creates two new blank rows (16,17) inheriting nearby formatting.
$ws.Range("16:17").Insert()

# --- Update row 13: date + text changes ---
$ws.Range("D13").Value = 44264
$ws.Range("E13").Value = "Finalization de la grille"

# --- Update row 14 ---
$ws.Range("D14").Value = 44267
$ws.Range("E14").Value = "Fin de la première version"

# --- Row 17 (the "log" entry, now at a new spot) ---
$ws.Range("D17").Value = 44287
$ws.Range("E17").Value = "Ajout d'un système de log"

# --- Row 18 (previously row 16, shifted down by the insert) ---
$ws.Range("E18").Value = "Ajout des fichier stoquant des grilles"

# --- Fill the newly inserted row 16 ---
$ws.Range("D16").Value = 44285
$ws.Range("E16").Value = "Ajout d'un système d'authentification"

# --- Update row 15 ---
$ws.Range("D15").Value = 44272
$ws.Range("E15").Value = "Implémentation du score"

# --- New VIP value on row 18 ---
$ws.Range("F18").Value = "VWM"
$ws.Range("F12").Copy()
$ws.Range("F18").PasteSpecial(-4122)

# --- Row height fix-ups ---
# Row 13 no longer needs its old 2-line height -> autofit back to default.
$ws.Rows.Item(13).AutoFit()

# Rows 16 and 18 now hold longer text wrapping onto two lines.
$ws.Rows.Item(16).RowHeight = 28.8
$ws.Rows.Item(18).RowHeight = 28.8

# --- Sheet view: scroll back to top-left and move the selection ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F19").Select()
